$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: Task / Time(h) / Date / Time(start)
$ws.Range("B11").Value = "Texture Loading"
$ws.Range("C11").Value = 3

# Reuse the existing date format (style used by D10) for D11
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D11").Value = 45679

# Reuse the existing time format (style used by E10) for E11
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E11").Value = 0.83333333333333337

# Match the author's final selection
$ws.Range("C3:C11").Select()
